$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in header cell B1: "locacalizacion" -> "localizacion"
$ws.Range("B1").Value = "localizacion"

# Update selection to just B1 (active cell B1, selection B1)
$ws.Range("B1").Select()
